$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Metano) - update fraction and dependent calculated values
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.7111383119035047
$ws.Range("H2").Value = 9315.118065181918
$ws.Range("I2").Value = 6624.337336055311
$ws.Range("J2").Value = 0.4301448194744298
$ws.Range("L2").Value = 10530.32048949866

# Row 3 (Etano) - update fraction and dependent calculated values
$ws.Range("F3").Value = 0.3
$ws.Range("G3").Value = 0.8910542058256088
$ws.Range("H3").Value = 22541.70222848128
$ws.Range("I3").Value = 20085.87857715674
$ws.Range("J3").Value = 0.6498681911765604

# Row 4 - Isobutano -> Propano, with new property values
$ws.Range("A4").Value = "Propano"
$ws.Range("B4").Value = 44.097
$ws.Range("C4").Value = 616
$ws.Range("D4").Value = 666.0599999999999
$ws.Range("E4").Value = 0.1522
$ws.Range("F4").Value = 0.1
$ws.Range("G4").Value = 0.9928385045840092
$ws.Range("H4").Value = 37926.76834742798
$ws.Range("I4").Value = 37655.15596976453
$ws.Range("J4").Value = 0.9027563927797403
$ws.Range("L4").Value = 0.6361476736065229

# Row 5 - new component "Heptano Plus"
$ws.Range("A5").Value = "Heptano Plus"
$ws.Range("B5").Value = 237
$ws.Range("C5").Value = 563
$ws.Range("D5").Value = 916
$ws.Range("E5").Value = 0.52
$ws.Range("F5").Value = 0.1
$ws.Range("G5").Value = 1.339084988163872
$ws.Range("H5").Value = 78484.22924214284
$ws.Range("I5").Value = 105097.0531857655
$ws.Range("J5").Value = 1.358391672383659
$ws.Range("A5:J5").Style = "Normal"

# Row 6 - A value (Peng-Robinson calculation)
$ws.Range("L6").Value = 0.9233872578344503

# Row 8 - B value
$ws.Range("L8").Value = 0.4034974477320561

# Row 10 - mc value
$ws.Range("L10").Value = 45.1522

# Row 12 - z value
$ws.Range("L12").Value = 1.070030438770379

# Row 14 - densidad value
$ws.Range("L14").Value = 26.76490044919463
